# Replace the abbreviated fiscal-month codes in column B (rows 2-64) with
# their full month names, and give that column a plain (non-bold, default
# color) Calibri font so it no longer inherits the bold/white header look.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$monthMap = @{
    "OCT" = "October";
    "NOV" = "November";
    "DEC" = "December";
    "JAN" = "January";
    "FEB" = "February";
    "MAR" = "March";
    "APR" = "April";
    "MAY" = "May";
    "JUN" = "June";
    "JUL" = "July";
    "AUG" = "August";
    "SEP" = "September";
}

$lastRow = $ws.Cells.Item($ws.Rows.Count, 2).End(-4162).Row
if ($lastRow -lt 64) { $lastRow = 64 }

for ($r = 2; $r -le $lastRow; $r++) {
    $cell = $ws.Cells.Item($r, 2)
    $old = $cell.Value()
    if ($monthMap.ContainsKey($old)) {
        $cell.Value = $monthMap[$old]
    }
}

$rng = $ws.Range("B2:B$lastRow")
$rng.Font.Name = "Calibri"
